$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins, Losses, Ties in AC1:AE1, matching style of existing header cells
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows: team record (constant 70-92-0) applied to every player row
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 29).Value = 70   # AC = Wins
    $ws.Cells.Item($r, 30).Value = 92   # AD = Losses
    $ws.Cells.Item($r, 31).Value = 0    # AE = Ties
}
